# Insert a new data row at row 93 (pushes existing rows 93-127 down to 94-128)
# and populate it with the new "Granada" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(93).Insert()

$ws.Cells.Item(93, 1).Value  = 10
$ws.Cells.Item(93, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(93, 3).Value  = 'La Araucanía'
$ws.Cells.Item(93, 4).Value  = '2022-06-14'
$ws.Cells.Item(93, 5).Value  = 9
$ws.Cells.Item(93, 6).Value  = 'Fruta'
$ws.Cells.Item(93, 7).Value  = 100104
$ws.Cells.Item(93, 8).Value  = 'Frutos de pepita'
$ws.Cells.Item(93, 9).Value  = 100104001
$ws.Cells.Item(93, 10).Value = 'Granada'
$ws.Cells.Item(93, 11).Value = 'Wonderfull'
$ws.Cells.Item(93, 12).Value = 'Especial'
$ws.Cells.Item(93, 13).Value = 50
$ws.Cells.Item(93, 14).Value = 15000
$ws.Cells.Item(93, 15).Value = 15000
$ws.Cells.Item(93, 16).Value = 15000
$ws.Cells.Item(93, 17).Value = '$/bandeja 10 kilos empedrada'
$ws.Cells.Item(93, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(93, 19).Value = 1500
$ws.Cells.Item(93, 20).Value = 10
